# Auto-generated script applying crypto price/volume refresh changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'90.454.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.93%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.107.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.55%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'237.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +7.94%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'626.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.15%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +5.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -8.38%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'3.105.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.710"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.33%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.50%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'36.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'5.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.29%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000241"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.92%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.162.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.677.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.72%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.119.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.17%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0000211"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.78%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'446.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.98%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +7.98%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'8.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.51%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'6.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.13%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'89.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.268.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -6.10%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'27.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +14.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.197"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +27.94%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +5.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'506.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "'3.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.70%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.77%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.29%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +10.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'22.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0851"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.10%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.735"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -18.22%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +38.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.693"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +11.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'148.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.92%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'45.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.56%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'4.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +7.01%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.60%  "
$ws.Range("E51").Style = "Normal"
